$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values to match repulled data
$ws.Range("F2").Value = -7
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 4
$ws.Range("F8").Value = 4
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = -1
$ws.Range("F14").Value = 1
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = 6
$ws.Range("F18").Value = -1
$ws.Range("F19").Value = 0
$ws.Range("F21").Value = 1
$ws.Range("F22").Value = 0
